$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.500.72"
$ws.Range("E2").Value = "  +1.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.547.44"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'571.24"
$ws.Range("E5").Value = "  +2.48%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'150.66"
$ws.Range("E6").Value = "  +8.22%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.74%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.544.78"
$ws.Range("E9").Value = "  +4.59%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.04%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.76"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.04%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "'0.360"
$ws.Range("E13").Value = "  +3.12%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +8.09%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.002.82"
$ws.Range("E15").Value = "  +4.80%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "63.424.16"
$ws.Range("E16").Value = "  +1.63%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +1.95%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.555.59"
$ws.Range("E18").Value = "  +4.93%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'340.91"
$ws.Range("E20").Value = "  -1.92%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  +3.84%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'6.86"
$ws.Range("E22").Value = "  +0.60%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'66.16"
$ws.Range("E24").Value = "  +1.29%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -1.01%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +5.60%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'8.66"
$ws.Range("E27").Value = "  +5.60%  "

# Row 28 - SuiNetwork
$ws.Range("D28").Value = "'1.50"
$ws.Range("E28").Value = "  +11.54%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.07%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = "  +12.13%  "

# Row 31 - PEPE
$ws.Range("D31").Value = "0.0₃0834"
$ws.Range("E31").Value = "  +5.63%  "

# Row 32 - PancakeSwap
$ws.Range("D32").Value = "'1.89"
$ws.Range("E32").Value = "  +4.43%  "

# Row 33 - Monero
$ws.Range("D33").Value = "'178.22"
$ws.Range("E33").Value = "  +3.61%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  +8.05%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "'422.39"
$ws.Range("E35").Value = "  +11.72%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("D36").Value = "'0.407"
$ws.Range("E36").Value = "  +2.36%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "'19.21"
$ws.Range("E37").Value = "  +3.30%  "

# Row 38 - NEARProtocol
$ws.Range("D38").Value = "'4.47"
$ws.Range("E38").Value = "  -0.68%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "'1.79"
$ws.Range("E40").Value = "  +5.55%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - OKB
$ws.Range("D42").Value = "'39.73"
$ws.Range("E42").Value = "  +1.49%  "

# Row 43 - Aave
$ws.Range("D43").Value = "'154.31"
$ws.Range("E43").Value = "  +5.78%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'3.82"
$ws.Range("E44").Value = "  +4.01%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "'21.18"
$ws.Range("E45").Value = "  +1.56%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "'0.611"
$ws.Range("E46").Value = "  +3.49%  "

# Row 47 - Hedera
$ws.Range("D47").Value = "'0.0531"
$ws.Range("E47").Value = "  +2.08%  "

# Row 48 - was VeChain, now Stellar (rows 48/49 swap content)
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.0969"
$ws.Range("E48").Value = "  +1.56%  "

# Row 49 - was Stellar, now VeChain
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0239"
$ws.Range("E49").Value = "  +7.20%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'18.69"
$ws.Range("E50").Value = "  +4.24%  "

# Row 51 - dogwifhat
$ws.Range("E51").Value = "  +7.15%  "
